$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2253521126760563
$ws.Range("C2").Value = 0.5154929577464789
$ws.Range("J2").Value = 0.01690140845070422
$ws.Range("P2").Value = 0.1577464788732394
$ws.Range("S2").Value = 0.08450704225352113
# Row 3
$ws.Range("B3").Value = 0.01081081081081081
$ws.Range("C3").Value = 0.01081081081081081
$ws.Range("J3").Value = 0.04324324324324325
$ws.Range("P3").Value = 0.7297297297297297
$ws.Range("S3").Value = 0.2054054054054054
# Row 4
$ws.Range("J4").Value = 0.1020408163265306
$ws.Range("P4").Value = 0.673469387755102
$ws.Range("S4").Value = 0.2244897959183673
# Row 6
$ws.Range("B6").Value = 0.06806282722513089
$ws.Range("D6").Value = 0.01570680628272251
$ws.Range("F6").Value = 0.08376963350785341
$ws.Range("J6").Value = 0.1465968586387434
$ws.Range("O6").Value = 0.01570680628272251
$ws.Range("Q6").Value = 0.193717277486911
$ws.Range("R6").Value = 0.06806282722513089
$ws.Range("S6").Value = 0.4083769633507853
# Row 7
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.00966183574879227
$ws.Range("F7").Value = 0.05314009661835749
$ws.Range("J7").Value = 0.1835748792270532
$ws.Range("O7").Value = 0.02898550724637681
$ws.Range("Q7").Value = 0.1449275362318841
$ws.Range("R7").Value = 0.0821256038647343
$ws.Range("S7").Value = 0.3864734299516908
# Row 8
$ws.Range("B8").Value = 0.09977827050997783
$ws.Range("D8").Value = 0.02439024390243903
$ws.Range("F8").Value = 0.05764966740576496
$ws.Range("J8").Value = 0.1263858093126386
$ws.Range("O8").Value = 0.008869179600886918
$ws.Range("Q8").Value = 0.188470066518847
$ws.Range("R8").Value = 0.09090909090909091
$ws.Range("S8").Value = 0.4035476718403548
# Row 9
$ws.Range("B9").Value = 0.1059907834101382
$ws.Range("D9").Value = 0.01382488479262673
$ws.Range("F9").Value = 0.04608294930875576
$ws.Range("J9").Value = 0.1382488479262673
$ws.Range("O9").Value = 0.01382488479262673
$ws.Range("Q9").Value = 0.2165898617511521
$ws.Range("R9").Value = 0.07373271889400922
$ws.Range("S9").Value = 0.391705069124424
# Row 10
$ws.Range("B10").Value = 0.1176890156918688
$ws.Range("D10").Value = 0.02211126961483595
$ws.Range("F10").Value = 0.06847360912981455
$ws.Range("J10").Value = 0.1412268188302425
$ws.Range("O10").Value = 0.01283880171184023
$ws.Range("Q10").Value = 0.1982881597717546
$ws.Range("R10").Value = 0.07631954350927246
$ws.Range("S10").Value = 0.3630527817403709
# Row 11
$ws.Range("G11").Value = 0.1201413427561837
$ws.Range("J11").Value = 0.09187279151943463
$ws.Range("K11").Value = 0.166077738515901
$ws.Range("L11").Value = 0.6113074204946997
$ws.Range("S11").Value = 0.01060070671378092
# Row 12
$ws.Range("G12").Value = 0.8176795580110497
$ws.Range("J12").Value = 0.1325966850828729
$ws.Range("K12").Value = 0.005524861878453038
$ws.Range("L12").Value = 0.03867403314917127
$ws.Range("S12").Value = 0.005524861878453038
# Row 13
$ws.Range("G13").Value = 0.6086956521739131
$ws.Range("J13").Value = 0.2826086956521739
$ws.Range("S13").Value = 0.108695652173913
# Row 14
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("S14").Value = 0.3333333333333333
# Row 15
$ws.Range("F15").Value = 0.004587155963302753
$ws.Range("H15").Value = 0.1100917431192661
$ws.Range("I15").Value = 0.1009174311926606
$ws.Range("J15").Value = 0.3853211009174312
$ws.Range("K15").Value = 0.07798165137614679
$ws.Range("M15").Value = 0.02752293577981652
$ws.Range("N15").Value = 0.004587155963302753
$ws.Range("O15").Value = 0.06422018348623854
$ws.Range("S15").Value = 0.2247706422018349
# Row 16
$ws.Range("F16").Value = 0.004651162790697674
$ws.Range("H16").Value = 0.2186046511627907
$ws.Range("I16").Value = 0.07906976744186046
$ws.Range("J16").Value = 0.3720930232558139
$ws.Range("K16").Value = 0.1348837209302326
$ws.Range("M16").Value = 0.01395348837209302
$ws.Range("N16").Value = 0.004651162790697674
$ws.Range("O16").Value = 0.05581395348837209
$ws.Range("S16").Value = 0.1162790697674419
# Row 17
$ws.Range("F17").Value = 0.008403361344537815
$ws.Range("H17").Value = 0.1491596638655462
$ws.Range("I17").Value = 0.1218487394957983
$ws.Range("J17").Value = 0.4453781512605042
$ws.Range("K17").Value = 0.07773109243697479
$ws.Range("M17").Value = 0.01890756302521008
$ws.Range("O17").Value = 0.05672268907563025
$ws.Range("S17").Value = 0.1218487394957983
# Row 18
$ws.Range("H18").Value = 0.2020725388601036
$ws.Range("I18").Value = 0.06217616580310881
$ws.Range("J18").Value = 0.4455958549222798
$ws.Range("K18").Value = 0.1088082901554404
$ws.Range("M18").Value = 0.02072538860103627
$ws.Range("N18").Value = 0.005181347150259068
$ws.Range("O18").Value = 0.05181347150259067
$ws.Range("S18").Value = 0.1036269430051813
# Row 19
$ws.Range("F19").Value = 0.003725782414307004
$ws.Range("H19").Value = 0.2049180327868853
$ws.Range("I19").Value = 0.08122205663189269
$ws.Range("J19").Value = 0.3926974664679583
$ws.Range("K19").Value = 0.09538002980625931
$ws.Range("M19").Value = 0.01862891207153502
$ws.Range("O19").Value = 0.06706408345752608
$ws.Range("S19").Value = 0.1363636363636364
